$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 72.22221999999999
$ws.Range("J9").Value = 149.5
$ws.Range("L9").Value = 149.5
$ws.Range("N9").Value = -487.5

$ws.Range("H74").Value = 129237.5
$ws.Range("I74").Value = 5650
$ws.Range("K74").Value = 5650
$ws.Range("M74").Value = -4714

$ws.Range("H77").Value = 129237.5
$ws.Range("I77").Value = 5650
$ws.Range("K77").Value = 28250
$ws.Range("M77").Value = -23570

$ws.Range("H138").Value = 1374.3043
$ws.Range("I138").Value = 561.9
$ws.Range("K138").Value = 1685.7
$ws.Range("M138").Value = 3454.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1064.9474
$ws.Range("I2").Value = 1064.9474
$ws.Range("K2").Value = 1064.9474
$ws.Range("M2").Value = -951.9474

$ws.Range("H5").Value = 113.333336
$ws.Range("I5").Value = 75
$ws.Range("J5").Value = 190
$ws.Range("K5").Value = 75
$ws.Range("L5").Value = 190
$ws.Range("M5").Value = 37
$ws.Range("N5").Value = -414

$ws.Range("H102").Value = 3225
$ws.Range("I102").Value = 3200
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 3200
$ws.Range("L102").Value = 3250
$ws.Range("M102").Value = -1578
$ws.Range("N102").Value = -6494

$ws.Range("H116").Value = 1064.9474
$ws.Range("I116").Value = 1064.9474
$ws.Range("K116").Value = 1064.9474
$ws.Range("M116").Value = 1229.0526

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1064.9474
$ws.Range("I3").Value = 1064.9474
$ws.Range("K3").Value = 1064.9474
$ws.Range("M3").Value = -950.9474

$ws.Range("H4").Value = 113.333336
$ws.Range("I4").Value = 75
$ws.Range("J4").Value = 190
$ws.Range("K4").Value = 75
$ws.Range("L4").Value = 190
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = -420

$ws.Range("H99").Value = 1485.5714
$ws.Range("I99").Value = 1483.1666
$ws.Range("K99").Value = 1483.1666
$ws.Range("M99").Value = 14.83339999999998

$ws.Range("H103").Value = 12950
$ws.Range("J103").Value = 12950
$ws.Range("L103").Value = 12950
$ws.Range("N103").Value = -15294

$ws.Range("H134").Value = 7613.222
$ws.Range("J134").Value = 6295.4443
$ws.Range("L134").Value = 18886.3329
$ws.Range("N134").Value = -23956.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 3500
$ws.Range("K31").Value = 500
$ws.Range("L31").Value = 3500
$ws.Range("M31").Value = -205
$ws.Range("N31").Value = -4090

$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 3500
$ws.Range("K34").Value = 500
$ws.Range("L34").Value = 3500
$ws.Range("M34").Value = -298
$ws.Range("N34").Value = -3904

$ws.Range("H58").Value = 10009.2
$ws.Range("I58").Value = 6342.5713
$ws.Range("J58").Value = 18564.666
$ws.Range("K58").Value = 6342.5713
$ws.Range("L58").Value = 18564.666
$ws.Range("M58").Value = -6139.5713
$ws.Range("N58").Value = -18970.666

$ws.Range("H136").Value = 10009.2
$ws.Range("I136").Value = 6342.5713
$ws.Range("J136").Value = 18564.666
$ws.Range("K136").Value = 19027.7139
$ws.Range("L136").Value = 55693.99800000001
$ws.Range("M136").Value = -16477.7139
$ws.Range("N136").Value = -60793.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 100
$ws.Range("I51").Value = 100
$ws.Range("K51").Value = 300
$ws.Range("M51").Value = 160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 182.77777
$ws.Range("I107").Value = 194.5
$ws.Range("K107").Value = 194.5
$ws.Range("M107").Value = 1725.5

$ws.Range("H126").Value = 1862.6
$ws.Range("I126").Value = 1578.25
$ws.Range("K126").Value = 4734.75
$ws.Range("M126").Value = -2264.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2844.4443
$ws.Range("I46").Value = 2228.5715
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 2228.5715
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -2040.5715
$ws.Range("N46").Value = -5376

$ws.Range("H61").Value = 4148.3
$ws.Range("I61").Value = 3218.8572
$ws.Range("J61").Value = 6317
$ws.Range("K61").Value = 3218.8572
$ws.Range("L61").Value = 6317
$ws.Range("M61").Value = -3016.8572
$ws.Range("N61").Value = -6721

$ws.Range("H68").Value = 4025
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 6550
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 6550
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -8048

$ws.Range("H71").Value = 4025
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 6550
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 32750
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -40238

$ws.Range("H104").Value = 35000
$ws.Range("J104").Value = 35000
$ws.Range("L104").Value = 35000
$ws.Range("N104").Value = -41988

$ws.Range("H113").Value = 4148.3
$ws.Range("I113").Value = 3218.8572
$ws.Range("J113").Value = 6317
$ws.Range("K113").Value = 3218.8572
$ws.Range("L113").Value = 6317
$ws.Range("M113").Value = -1048.8572
$ws.Range("N113").Value = -10657

$ws.Range("H122").Value = 3263.818
$ws.Range("I122").Value = 2500.4285
$ws.Range("K122").Value = 7501.2855
$ws.Range("M122").Value = -5051.2855

$ws.Range("H132").Value = 10849.75
$ws.Range("I132").Value = 10849.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 32549.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -30019.25
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 500
$ws.Range("J14").Value = 500
$ws.Range("L14").Value = 500
$ws.Range("N14").Value = -836

$ws.Range("H80").Value = 79998.336
$ws.Range("J80").Value = 79998.336
$ws.Range("L80").Value = 79998.336
$ws.Range("N80").Value = -81994.336

$ws.Range("H83").Value = 79998.336
$ws.Range("J83").Value = 79998.336
$ws.Range("L83").Value = 239995.008
$ws.Range("N83").Value = -249979.008

$ws.Range("H104").Value = 11833.333
$ws.Range("J104").Value = 11833.333
$ws.Range("L104").Value = 11833.333
$ws.Range("N104").Value = -18821.333

$ws.Range("H132").Value = 1062.25
$ws.Range("I132").Value = 1071.1428
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3213.4284
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -683.4284000000002
$ws.Range("N132").Value = -8060
